# Contabilità workbook update
# - Add a new "balance" data row (3000 / 500 / 2500)
# - Add a message row showing the label "Saldo" with its value (3500)
# - Rename the header of column C from "Saldo" to "Differenza"
# - Move the selection off the old (now removed) filtered balance range
#   and onto a single cell, since the data rows are now dynamic

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row
$ws.Range("A4").Value = 3000
$ws.Range("B4").Value = 500
$ws.Range("C4").Value = 2500

# New "Saldo" message row with its value
$ws.Range("B5").Value = "Saldo"
$ws.Range("C5").Value = 3500

# Column header C1: "Saldo" -> "Differenza"
$ws.Range("C1").Value = "Differenza"

# Reset selection to a single cell below the data
$ws.Range("D6").Select() | Out-Null
